# Updated solution for Tutorial 6
# Dates change from DD/MM/YYYY to DD-MM-YYYY text format, and some
# attendance flag columns (D/E/G/H) are corrected for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the date cells as literal text (matching the original inline-string
# storage) instead of letting Excel auto-parse ambiguous day<=12 strings
# like "01-08-2022" into real date serials.
$cell = $ws.Range("A3")
$cell.NumberFormat = "@"
$cell.Value = "28-07-2022"
$cell.Style = "Normal"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$cell = $ws.Range("A4")
$cell.NumberFormat = "@"
$cell.Value = "01-08-2022"
$cell.Style = "Normal"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$cell = $ws.Range("A5")
$cell.NumberFormat = "@"
$cell.Value = "04-08-2022"
$cell.Style = "Normal"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$cell = $ws.Range("A6")
$cell.NumberFormat = "@"
$cell.Value = "08-08-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A7")
$cell.NumberFormat = "@"
$cell.Value = "11-08-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A8")
$cell.NumberFormat = "@"
$cell.Value = "15-08-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A9")
$cell.NumberFormat = "@"
$cell.Value = "18-08-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A10")
$cell.NumberFormat = "@"
$cell.Value = "22-08-2022"
$cell.Style = "Normal"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

$cell = $ws.Range("A11")
$cell.NumberFormat = "@"
$cell.Value = "25-08-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A12")
$cell.NumberFormat = "@"
$cell.Value = "29-08-2022"
$cell.Style = "Normal"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

$cell = $ws.Range("A13")
$cell.NumberFormat = "@"
$cell.Value = "01-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A14")
$cell.NumberFormat = "@"
$cell.Value = "05-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A15")
$cell.NumberFormat = "@"
$cell.Value = "08-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A16")
$cell.NumberFormat = "@"
$cell.Value = "12-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A17")
$cell.NumberFormat = "@"
$cell.Value = "15-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A18")
$cell.NumberFormat = "@"
$cell.Value = "19-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A19")
$cell.NumberFormat = "@"
$cell.Value = "22-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A20")
$cell.NumberFormat = "@"
$cell.Value = "26-09-2022"
$cell.Style = "Normal"

$cell = $ws.Range("A21")
$cell.NumberFormat = "@"
$cell.Value = "29-09-2022"
$cell.Style = "Normal"
